$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "rxxx"
$ws.Range("B8").Value = "greg"
$ws.Range("C8").Value = "temp profile off on setup"
$ws.Range("D8").Value = "2025-09-30 13:25:43"
